$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: write new/updated cell values (formats fixed up in step 2) ---
$ws.Range("C7").Value = "-Code UI"
$ws.Range("B8").Value = "Tuần 9"
$ws.Range("C8").Value = "-Create database with MySQL`n-Code Front-end with React"
$ws.Range("D8").Value = "- xong CSDL`n'-xong UI and Admin page"
$ws.Range("B9").Value = "Tuần 10"
$ws.Range("C9").Value = "-Code Back-end API full`n-Code Front-end full page"
$ws.Range("B10").Value = "Tuần 11"
$ws.Range("C10").Value = "-Database hoàn thiện"
$ws.Range("B11").Value = "Tuần 12"
$ws.Range("C11").Value = "-Code hoàn thiện"
$ws.Range("B12").Value = "Tuần 13"
$ws.Range("C12").Value = "-Báo cáo hoàn thiện"
$ws.Range("B13").Value = "Tuần 14"
$ws.Range("C13").Value = "-Xong!"

# --- Step 2: copy cell formatting from existing reference cells so the shared-style indices match exactly ---
# C7 lost its quotePrefix style when the value was reassigned -> restore from D7 (same row/style family)
$ws.Range("D7").Copy()
$ws.Range("C7").PasteSpecial(-4122)

# Wrap-text + border style (like C5/C6/D5/D6) -> C8, D8, C9
$ws.Range("C6").Copy()
$ws.Range("C8").PasteSpecial(-4122)
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)

# Border-only (quotePrefix) style (like C7) -> C10, C11, C12, C13
$ws.Range("C7").Copy()
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("C13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Step 3: row heights for the two newly wrapped rows ---
$ws.Rows.Item(8).RowHeight = 37.5
$ws.Rows.Item(9).RowHeight = 37.5

# --- Step 4: view state (zoom + selection) ---
$ws.Application.ActiveWindow.Zoom = 85
[void]$ws.Range("D9").Select()

